# Append the new "Daily APR" data row (id=5) to the active worksheet,
# matching the upstream data export that produced this workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the next empty row right after the current data (row 5 -> row 6).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$newRow = $lastRow + 1

$ws.Cells.Item($newRow, 1).Value = 5
$ws.Cells.Item($newRow, 2).Value = "2025-09-01T22:16"
$ws.Cells.Item($newRow, 3).Value = 1.6714039663513292
